# Added try/except around the Excel write logic (per commit message).
# The underlying data changes are a re-run of the same backward-elimination
# report a few days later: the absPath drive letter casing, the revision
# GUID, and every embedded "Date:"/"Time:" stamp inside each sheet's OLS
# summary text (B2) move from "Sun, 29 Dec 2019 16:11:14" to
# "Wed, 01 Jan 2020 23:18:52" (sheets 1-11) / "23:18:53" (sheets 12-16).

try {
    $wb = $excel.ActiveWorkbook

    foreach ($ws in $wb.Worksheets) {
        if ($ws.Index -le 11) {
            $newTime = "23:18:52"
        } else {
            $newTime = "23:18:53"
        }

        $cell = $ws.Range("B2")
        $text = $cell.Value()

        if ($text) {
            $text = $text.Replace("Sun, 29 Dec 2019", "Wed, 01 Jan 2020")
            $text = $text.Replace("16:11:14", $newTime)
            $cell.Value = $text
        }
    }
}
catch {
    Write-Output $_
}
